# Add a new elabe poll (fieldwork 11/11) made of three rows (105-107),
# matching the "add elabe poll (11/11)" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 105 ---------------------------------------------------------
$ws.Range("A105").Value  = 31
$ws.Range("B105").Value  = 2021
$ws.Range("C105").Value  = 11
$ws.Range("D105").Value  = 11
$ws.Range("E105").Value  = 11
$ws.Range("F105").Value  = "elabe"
$ws.Range("G105").Value  = "online"
$ws.Range("H105").Value  = "partially"
$ws.Range("I105").Value  = 933
$ws.Range("J105").Value  = 2
$ws.Range("K105").Value  = 1
$ws.Range("L105").Value  = 8
$ws.Range("M105").Value  = 2
$ws.Range("N105").Value  = 2
$ws.Range("O105").Value  = 7
$ws.Range("P105").Value  = 5
$ws.Range("Q105").Value  = 25
$ws.Range("T105").Value  = 13
$ws.Range("U105").Value  = 1
$ws.Range("V105").Value  = 2
$ws.Range("W105").Value  = 17
$ws.Range("X105").Value  = 14
$ws.Range("Y105").Value  = 0.5
$ws.Range("Z105").Value  = "T_0.5"
$ws.Range("AA105").Value = 0.5

# --- Row 106 ---------------------------------------------------------
$ws.Range("A106").Value  = 31
$ws.Range("B106").Value  = 2021
$ws.Range("C106").Value  = 11
$ws.Range("D106").Value  = 11
$ws.Range("E106").Value  = 11
$ws.Range("F106").Value  = "elabe"
$ws.Range("G106").Value  = "online"
$ws.Range("H106").Value  = "partially"
$ws.Range("I106").Value  = 919
$ws.Range("J106").Value  = 2
$ws.Range("K106").Value  = 1
$ws.Range("L106").Value  = 8
$ws.Range("M106").Value  = 2
$ws.Range("N106").Value  = 2
$ws.Range("O106").Value  = 7
$ws.Range("P106").Value  = 5
$ws.Range("Q106").Value  = 27
$ws.Range("R106").Value  = 11
$ws.Range("U106").Value  = 0.5
$ws.Range("V106").Value  = 2
$ws.Range("W106").Value  = 18
$ws.Range("X106").Value  = 14
$ws.Range("Y106").Value  = "T_0.5"
$ws.Range("Z106").Value  = "T_0.5"
$ws.Range("AA106").Value = 0.5

# --- Row 107 ---------------------------------------------------------
$ws.Range("A107").Value  = 31
$ws.Range("B107").Value  = 2021
$ws.Range("C107").Value  = 11
$ws.Range("D107").Value  = 11
$ws.Range("E107").Value  = 11
$ws.Range("F107").Value  = "elabe"
$ws.Range("G107").Value  = "online"
$ws.Range("H107").Value  = "partially"
$ws.Range("I107").Value  = 933
$ws.Range("J107").Value  = 2
$ws.Range("K107").Value  = 1
$ws.Range("L107").Value  = 8
$ws.Range("M107").Value  = 2
$ws.Range("N107").Value  = 2
$ws.Range("O107").Value  = 8
$ws.Range("P107").Value  = 5
$ws.Range("Q107").Value  = 28
$ws.Range("S107").Value  = 9
$ws.Range("U107").Value  = 1
$ws.Range("V107").Value  = 2
$ws.Range("W107").Value  = 18
$ws.Range("X107").Value  = 13
$ws.Range("Y107").Value  = 0.5
$ws.Range("Z107").Value  = "T_0.5"
$ws.Range("AA107").Value = 0.5

# The "T_0.5" text markers reuse the workbook's existing explicit-black
# font style (rather than the theme-color default font) so the saved
# cell style index matches the other "T_0.5"/"T_1" cells in the sheet.
$ws.Range("Z105").Font.Color = 0
$ws.Range("Y106").Font.Color = 0
$ws.Range("Z106").Font.Color = 0
$ws.Range("Z107").Font.Color = 0

# Restore the sheet selection to just past the new data, as in the
# saved workbook (selection moves from U104 to AA108).
[void]$ws.Range("AA108").Select()
